# Update the 15.a.1 indicator text (cell B4 on the single worksheet) to the
# revised wording. The old shared-string entry is dropped and a new one is
# appended, and the selection moves from B2 to B4 (matching the author's
# workbook view when they saved after editing this cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = "15.a.1 (a) Official development assistance on conservation and sustainable use of biodiversity; and (b) revenue generated and finance mobilized from biodiversity-relevant economic instruments"

# Move the active selection to B4 (also clears the previous topLeftCell="A2"
# scroll-freeze that pointed the view at row 2).
$ws.Range("B4").Select()
